# Burndown chart / Backlog update:
#  - "Upload sprints to github (Todo mundo)"  -> "Upload sprints to github (Todo agente)"
#  - "Take metrics (João)"                    -> "Take metrics (Toda agente)"
# (the "Implement feature 2 (...)" row text is unchanged)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = "Upload sprints to github (Todo agente)"
$ws.Range("C9").Value = "Take metrics (Toda agente)"
